$wb = $excel.ActiveWorkbook
$demo = $wb.Worksheets.Item("Demo")
$demo.Name = "Expenses"
$demo.Range("C3").Value = 1000

$income = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $demo)
$income.Name = "Income"

$expRange = $demo.Range("A1:C4")
$expRange.Copy()
$income.Range("A1:C4").PasteSpecial(-4122)

$income.Range("A1").Value = "Row no"
$income.Range("B1").Value = "Description"
$income.Range("C1").Value = "Amount"

$income.Range("A2").Value = 3
$income.Range("B2").Value = "Salary"
$income.Range("C2").Value = 13000

$income.Range("A3").Value = "N/A"
$income.Range("B3").Value = "Pension"
$income.Range("C3").Value = 700

$income.Range("A4").Value = "N/A"
$income.Range("B4").Value = "Labour"
$income.Range("C4").Value = 5000

$income.Columns("A:C").AutoFit()
$income.Range("C3").Select()

$demo.Activate()
$demo.Range("A2").Select()
